# Fix bug: price list values were divided by an incorrect (exceeded) request
# factor from Google Drive; correct the date and recalculate the prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Date bumped by one day
$ws.Range("A1").Value = 45311

# Recalculated prices (column D)
$ws.Range("D33").Value = 3287.585
$ws.Range("D34").Value = 2470.096
$ws.Range("D35").Value = 2131.561
$ws.Range("D36").Value = 1914.227
$ws.Range("D37").Value = 1914.227
$ws.Range("D38").Value = 1614.548
$ws.Range("D39").Value = 1614.548
$ws.Range("D40").Value = 1614.548
$ws.Range("D41").Value = 1614.548
$ws.Range("D42").Value = 1614.548
$ws.Range("D43").Value = 1614.548
$ws.Range("D44").Value = 1614.548
$ws.Range("D45").Value = 1926.763
$ws.Range("D46").Value = 1926.763
$ws.Range("D47").Value = 1926.763
$ws.Range("D51").Value = 2808.644
$ws.Range("D52").Value = 2553.696
$ws.Range("D53").Value = 2131.561
$ws.Range("D54").Value = 2131.561
